$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 211
$ws.Cells.Item(2, 2).Value = 236
$ws.Cells.Item(2, 3).Value = 191
$ws.Cells.Item(2, 4).Value = 232
$ws.Cells.Item(2, 5).Value = 195
$ws.Cells.Item(2, 6).Value = 204
$ws.Cells.Item(2, 7).Value = 216
$ws.Cells.Item(2, 8).Value = 207

# Row 3
$ws.Cells.Item(3, 1).Value = 239
$ws.Cells.Item(3, 2).Value = 218
$ws.Cells.Item(3, 3).Value = 211
$ws.Cells.Item(3, 4).Value = 218
$ws.Cells.Item(3, 5).Value = 211
$ws.Cells.Item(3, 6).Value = 195
$ws.Cells.Item(3, 7).Value = 239
$ws.Cells.Item(3, 8).Value = 195

# Row 4
$ws.Cells.Item(4, 1).Value = 184
$ws.Cells.Item(4, 2).Value = 278
$ws.Cells.Item(4, 3).Value = 158
$ws.Cells.Item(4, 4).Value = 278
$ws.Cells.Item(4, 5).Value = 158
$ws.Cells.Item(4, 6).Value = 241
$ws.Cells.Item(4, 7).Value = 184
$ws.Cells.Item(4, 8).Value = 241

# Row 5
$ws.Cells.Item(5, 1).Value = 165
$ws.Cells.Item(5, 2).Value = 247
$ws.Cells.Item(5, 3).Value = 155
$ws.Cells.Item(5, 4).Value = 202
$ws.Cells.Item(5, 5).Value = 189
$ws.Cells.Item(5, 6).Value = 194
$ws.Cells.Item(5, 7).Value = 199
$ws.Cells.Item(5, 8).Value = 240

# Row 6
$ws.Cells.Item(6, 1).Value = 213
$ws.Cells.Item(6, 2).Value = 297
$ws.Cells.Item(6, 3).Value = 174
$ws.Cells.Item(6, 4).Value = 271
$ws.Cells.Item(6, 5).Value = 196
$ws.Cells.Item(6, 6).Value = 238
$ws.Cells.Item(6, 7).Value = 235
$ws.Cells.Item(6, 8).Value = 264

# Row 7
$ws.Cells.Item(7, 1).Value = 199
$ws.Cells.Item(7, 2).Value = 206
$ws.Cells.Item(7, 3).Value = 181
$ws.Cells.Item(7, 4).Value = 158
$ws.Cells.Item(7, 5).Value = 227
$ws.Cells.Item(7, 6).Value = 141
$ws.Cells.Item(7, 7).Value = 245
$ws.Cells.Item(7, 8).Value = 189

# Row 8
$ws.Cells.Item(8, 1).Value = 227
$ws.Cells.Item(8, 2).Value = 281
$ws.Cells.Item(8, 3).Value = 211
$ws.Cells.Item(8, 4).Value = 233
$ws.Cells.Item(8, 5).Value = 277
$ws.Cells.Item(8, 6).Value = 210
$ws.Cells.Item(8, 7).Value = 294
$ws.Cells.Item(8, 8).Value = 258

# Row 9
$ws.Cells.Item(9, 1).Value = 164
$ws.Cells.Item(9, 2).Value = 203
$ws.Cells.Item(9, 3).Value = 94
$ws.Cells.Item(9, 4).Value = 199
$ws.Cells.Item(9, 5).Value = 96
$ws.Cells.Item(9, 6).Value = 148
$ws.Cells.Item(9, 7).Value = 167
$ws.Cells.Item(9, 8).Value = 151

# Row 10
$ws.Cells.Item(10, 1).Value = 103
$ws.Cells.Item(10, 2).Value = 293
$ws.Cells.Item(10, 3).Value = 99
$ws.Cells.Item(10, 4).Value = 241
$ws.Cells.Item(10, 5).Value = 166
$ws.Cells.Item(10, 6).Value = 235
$ws.Cells.Item(10, 7).Value = 170
$ws.Cells.Item(10, 8).Value = 288

